$wb = $excel.ActiveWorkbook

# Build the target header format (bold font, thin box border, centered
# horizontally and aligned to the top) once on a scratch cell so the
# style table stays compact, then stamp that exact format onto every
# sheet's A1/B1 header cells via copy / paste-special (formats only).
$templateWs = $wb.Worksheets.Item(1)
$tmpl = $templateWs.Range("Z1")
$tmpl.Font.Bold = $true
$tmpl.HorizontalAlignment = -4108   # xlCenter
$tmpl.VerticalAlignment = -4160     # xlTop
$tmpl.Borders.LineStyle = 1         # xlContinuous (thin by default)
$tmpl.Copy()

foreach ($ws in $wb.Worksheets) {
    $a1 = $ws.Range("A1")
    $a1.Value = "Input Sheet"
    $b1 = $ws.Range("B1")
    $b1.Value = "Value"

    $a1.PasteSpecial(-4122)  # xlPasteFormats
    $b1.PasteSpecial(-4122)  # xlPasteFormats
}

$tmpl.Clear()
$excel.CutCopyMode = $false
